# Apply "upload orgincode detecting result" edit:
#  1. Add two new metric columns (F, G) of data for rows 14-16.
#  2. Add a summary row 19 (cols B:G) averaging rows 14-16 with a shared formula.
#  3. Re-position/resize the chart (graphicFrame) on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New data in columns F/G for rows 14-16 ---------------------------
$ws.Range("F14").Value = 0.69099999999999995
$ws.Range("G14").Value = 0.53100000000000003

$ws.Range("F15").Value = 0.69399999999999995
$ws.Range("G15").Value = 0.51700000000000002

$ws.Range("F16").Value = 0.68899999999999995
$ws.Range("G16").Value = 0.52400000000000002

# --- 2. Summary row 19: average (SUM/3) of rows 14:16 for each column ----
$ws.Range("B19").Formula = "=SUM(B14:B16)/3"
$ws.Range("C19:G19").Formula = "=SUM(C14:C16)/3"

# --- 3. Move & resize the existing chart ----------------------------------
$co = $ws.ChartObjects(1)
$co.Left = 633.8750787401575
$co.Top = 89.18559055118111
$co.Width = 385.0625
$co.Height = 216.10716535433073
